$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = "'-0.883"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(2,6).Value = -1.471
$ws.Cells.Item(2,7).Value = -0.29

# Row 3
$ws.Cells.Item(3,5).Value = "'-0.283"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,6).Value = -1.194
$ws.Cells.Item(3,7).Value = 0.606

# Row 4
$ws.Cells.Item(4,5).Value = "'0.009"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(4,6).Value = -0.644
$ws.Cells.Item(4,7).Value = 0.677

# Row 5
$ws.Cells.Item(5,5).Value = "'0.204"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(5,6).Value = -0.671
$ws.Cells.Item(5,7).Value = 1.059

# Row 6
$ws.Cells.Item(6,5).Value = "'98%"
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,5).Value = "'2%"
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,5).Value = "'-0.568"
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(8,6).Value = -1.72
$ws.Cells.Item(8,7).Value = 0.52

# Row 9
$ws.Cells.Item(9,5).Value = "'0.263"
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(9,6).Value = -1.181
$ws.Cells.Item(9,7).Value = 1.766

# Row 10
$ws.Cells.Item(10,5).Value = "'-0.548"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(10,6).Value = -1.802
$ws.Cells.Item(10,7).Value = 0.643

# Row 11
$ws.Cells.Item(11,5).Value = "'1.12"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(11,6).Value = -0.576
$ws.Cells.Item(11,7).Value = 2.873

# Row 14
$ws.Cells.Item(14,5).Value = "'1.542"
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(14,6).Value = 1.102
$ws.Cells.Item(14,7).Value = 1.985

# Row 15
$ws.Cells.Item(15,5).Value = "'0.289"
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(15,6).Value = -0.308
$ws.Cells.Item(15,7).Value = 0.892

# Row 16
$ws.Cells.Item(16,5).Value = "'0.595"
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(16,6).Value = 0.218
$ws.Cells.Item(16,7).Value = 0.968

# Row 17
$ws.Cells.Item(17,5).Value = "'0.229"
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(17,6).Value = -0.175
$ws.Cells.Item(17,7).Value = 0.625

# Row 18
$ws.Cells.Item(18,5).Value = "'-0.653"
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(18,6).Value = -1.151
$ws.Cells.Item(18,7).Value = -0.152

# Row 19
$ws.Cells.Item(19,5).Value = "'-0.029"
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(19,6).Value = -0.547
$ws.Cells.Item(19,7).Value = 0.508

# Row 24
$ws.Cells.Item(24,5).Value = "'1.637"
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(24,6).Value = 1.368
$ws.Cells.Item(24,7).Value = 1.893

# Row 25
$ws.Cells.Item(25,5).Value = "'-0.4"
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(25,6).Value = -0.772
$ws.Cells.Item(25,7).Value = -0.026

# Row 26
$ws.Cells.Item(26,5).Value = "'0.494"
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(26,6).Value = 0.2
$ws.Cells.Item(26,7).Value = 0.788

# Row 27
$ws.Cells.Item(27,6).Value = -0.208
$ws.Cells.Item(27,7).Value = 0.332

# Row 28
$ws.Cells.Item(28,5).Value = "'0.051"
$ws.Cells.Item(28,5).Style = "Normal"
$ws.Cells.Item(28,6).Value = -0.355
$ws.Cells.Item(28,7).Value = 0.453

# Row 29
$ws.Cells.Item(29,5).Value = "'0.201"
$ws.Cells.Item(29,5).Style = "Normal"
$ws.Cells.Item(29,6).Value = -0.183
$ws.Cells.Item(29,7).Value = 0.58

# Row 34
$ws.Cells.Item(34,5).Value = "'0.508"
$ws.Cells.Item(34,5).Style = "Normal"
$ws.Cells.Item(34,6).Value = 0.209
$ws.Cells.Item(34,7).Value = 0.8070000000000001

# Row 35
$ws.Cells.Item(35,5).Value = "'-0.09"
$ws.Cells.Item(35,5).Style = "Normal"
$ws.Cells.Item(35,6).Value = -0.524

# Row 36
$ws.Cells.Item(36,5).Value = "'0.138"
$ws.Cells.Item(36,5).Style = "Normal"
$ws.Cells.Item(36,6).Value = -0.107
$ws.Cells.Item(36,7).Value = 0.385

# Row 37
$ws.Cells.Item(37,5).Value = "'0.086"
$ws.Cells.Item(37,5).Style = "Normal"
$ws.Cells.Item(37,6).Value = -0.171
$ws.Cells.Item(37,7).Value = 0.347

# Row 38
$ws.Cells.Item(38,5).Value = "'-0.145"
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(38,6).Value = -0.477
$ws.Cells.Item(38,7).Value = 0.206

# Row 39
$ws.Cells.Item(39,5).Value = "'0.001"
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(39,6).Value = -0.357
$ws.Cells.Item(39,7).Value = 0.356
